# Apply the final-results update described in the commit "This should be it!
# all final results included."
#
# 1) Two worksheets ("100 Australian species" x CO1 and "Lutjanidae" x CO1)
#    had placeholder zero values for MMSeqs2_100 / MMSeqs2_97 (rows 8 & 9,
#    columns D-H) that are now filled in with the real computed metrics.
# 2) The "Rottnest" query/site label used across the three Rottnest-based
#    result sheets (12S, 16S, CO1) has been renamed to "Wadjemup" (the
#    island's indigenous name), and that sheet's MMSeqs2_100 / MMSeqs2_97
#    placeholder zeros (rows 8 & 9, columns D-H) are filled in as well.

$wb = $excel.ActiveWorkbook

# --- Sheet3: "100 Australian species" / CO1 -------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("D8").Value = 0.9722222222222222
$ws3.Range("E8").Value = 0.5
$ws3.Range("F8").Value = 0.660377358490566
$ws3.Range("G8").Value = 0.8177570093457943
$ws3.Range("H8").Value = 0.6363636363636364

$ws3.Range("D9").Value = 0.9830508474576272
$ws3.Range("E9").Value = 0.8285714285714286
$ws3.Range("F9").Value = 0.8992248062015504
$ws3.Range("G9").Value = 0.9477124183006537
$ws3.Range("H9").Value = 0.8686868686868687

# --- Sheet6: "Lutjanidae" / CO1 --------------------------------------------
$ws6 = $wb.Worksheets.Item("Sheet6")

$ws6.Range("D8").Value = 1
$ws6.Range("E8").Value = 0.44
$ws6.Range("F8").Value = 0.6111111111111112
$ws6.Range("G8").Value = 0.7971014492753624
$ws6.Range("H8").Value = 0.4814814814814815

$ws6.Range("D9").Value = 1
$ws6.Range("E9").Value = 0.92
$ws6.Range("F9").Value = 0.9583333333333334
$ws6.Range("G9").Value = 0.9829059829059831
$ws6.Range("H9").Value = 0.9259259259259259

# --- Sheets 7, 8, 9: rename "Rottnest" -> "Wadjemup" -----------------------
$rottnestSheets = @("Sheet7", "Sheet8", "Sheet9")
foreach ($sheetName in $rottnestSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 13; $row++) {
        $ws.Range("B$row").Value = "Wadjemup"
    }
}

# --- Sheet9: "Wadjemup" / CO1 zero-value fill-in ---------------------------
$ws9 = $wb.Worksheets.Item("Sheet9")

$ws9.Range("D8").Value = 0.9830508474576272
$ws9.Range("E8").Value = 0.5087719298245614
$ws9.Range("F8").Value = 0.6705202312138728
$ws9.Range("G8").Value = 0.8285714285714286
$ws9.Range("H8").Value = 0.5128205128205128

$ws9.Range("D9").Value = 0.979381443298969
$ws9.Range("E9").Value = 0.8407079646017699
$ws9.Range("F9").Value = 0.9047619047619049
$ws9.Range("G9").Value = 0.9481037924151696
$ws9.Range("H9").Value = 0.8290598290598291
